# "terceira parte 20 de abril de 2025"
# Adds a third column "precisao" = "100.0%" and resets the A2/B2 sample
# values (n_moscas / Pontos row) to 1 and 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 = "precisao", matching the bold/bordered header style
# already used by A1/B1.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "precisao"

# Row 2 data.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1

# C2 = "100.0%" as literal text (not an auto-converted percentage number).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "100.0%"
$ws.Range("C2").ClearFormats()
